# Inserts a new 9-row block of Pimiento price data (fecha 2023-04-05 / serial 45021)
# for "Comercializadora del Agro de Limarí" above the existing rows, shifting all
# subsequent rows down by 9 (dimension grows from R1282 to R1291).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 blank rows at the top of the existing Pimiento block (row 1208..1216)
$ws.Range("A1208:R1216").Insert()

# Columns: H (Variedad), I (Calidad), J (Volumen), K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$data = @(
    @("Cuatro cascos rojo",  "Primera", 400,  9000,  10000, 9500, 528),
    @("Cuatro cascos rojo",  "Segunda", 300,  7000,  8000,  7500, 417),
    @("Cuatro cascos rojo",  "Tercera", 200,  5000,  6000,  5500, 306),
    @("Cuatro cascos verde", "Primera", 700,  7000,  8000,  7500, 417),
    @("Cuatro cascos verde", "Segunda", 600,  5000,  6000,  5500, 306),
    @("Cuatro cascos verde", "Tercera", 400,  3000,  4000,  3500, 194),
    @("Morrón rojo",         "Primera", 1100, 8000,  9000,  8500, 472),
    @("Morrón rojo",         "Segunda", 600,  6000,  7000,  6500, 361),
    @("Morrón rojo",         "Tercera", 500,  4000,  5000,  4500, 250)
)

$r = 1208
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value2 = 2
    $ws.Cells.Item($r, 2).Value2 = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($r, 3).Value2 = "Coquimbo"
    $ws.Cells.Item($r, 4).Value2 = 45021
    $ws.Cells.Item($r, 5).Value2 = 4
    $ws.Cells.Item($r, 6).Value2 = 100112002
    $ws.Cells.Item($r, 7).Value2 = "Pimiento"
    $ws.Cells.Item($r, 8).Value2 = $row[0]
    $ws.Cells.Item($r, 9).Value2 = $row[1]
    $ws.Cells.Item($r, 10).Value2 = $row[2]
    $ws.Cells.Item($r, 11).Value2 = $row[3]
    $ws.Cells.Item($r, 12).Value2 = $row[4]
    $ws.Cells.Item($r, 13).Value2 = $row[5]
    $ws.Cells.Item($r, 14).Value2 = "`$/caja 18 kilos"
    $ws.Cells.Item($r, 15).Value2 = "Provincia de Limarí"
    $ws.Cells.Item($r, 16).Value2 = $row[6]
    $ws.Cells.Item($r, 17).Value2 = 18
    $ws.Cells.Item($r, 18).Value2 = "Hortaliza"
    $r = $r + 1
}
